# Update issue 260 - Update performance document.
# Adds a new "v1423" benchmark-run column to the "Sponza" and "ComplexMesh"
# worksheets (column L on Sponza, column K on ComplexMesh), matching the
# existing layout: 10 per-run values, an AVERAGE row, a VAR.S row, a
# 1-T.TEST row and two ratio rows, plus updated selection/dimension and
# conditional formatting ranges.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sponza (sheet2) - new column L, "v1423"
# ---------------------------------------------------------------------
$sponza = $wb.Worksheets.Item("Sponza")

# Extend formatting from the existing last column (K) into the new one (L)
# so the new column inherits the same header/body/summary styles.
$sponza.Range("K1:K16").Copy() | Out-Null
$sponza.Range("L1:L16").PasteSpecial(-4122) | Out-Null

$sponza.Range("L1").Value = "v1423"

$sponza.Range("L2").Value = 7242
$sponza.Range("L3").Value = 7229
$sponza.Range("L4").Value = 7226
$sponza.Range("L5").Value = 7254
$sponza.Range("L6").Value = 7204
$sponza.Range("L7").Value = 7226
$sponza.Range("L8").Value = 7212
$sponza.Range("L9").Value = 7262
$sponza.Range("L10").Value = 7229
$sponza.Range("L11").Value = 7264

$sponza.Range("L12").Formula = "=AVERAGE(L2:L11)"
$sponza.Range("L13").Formula = "=_xlfn.VAR.S(L2:L11)"
$sponza.Range("L14").Formula = "=1-_xlfn.T.TEST(K2:K11,L2:L11,2,3)"
$sponza.Range("L15").Formula = "=K12/L12"
$sponza.Range("L16").Formula = "=B12/L12"

$sponza.Range("B15:K16").FormatConditions.Item(1).ModifyAppliesToRange($sponza.Range("B15:L16")) | Out-Null

$sponza.Range("A1").Select() | Out-Null
$sponza.Range("L2").Select() | Out-Null

# ---------------------------------------------------------------------
# ComplexMesh (sheet3) - new column K, "v1423"
# ---------------------------------------------------------------------
$mesh = $wb.Worksheets.Item("ComplexMesh")

$mesh.Range("J1:J16").Copy() | Out-Null
$mesh.Range("K1:K16").PasteSpecial(-4122) | Out-Null

$mesh.Range("K1").Value = "v1423"

$mesh.Range("K2").Value = 5249
$mesh.Range("K3").Value = 5080
$mesh.Range("K4").Value = 5073
$mesh.Range("K5").Value = 5055
$mesh.Range("K6").Value = 5051
$mesh.Range("K7").Value = 5085
$mesh.Range("K8").Value = 5050
$mesh.Range("K9").Value = 5044
$mesh.Range("K10").Value = 5058
$mesh.Range("K11").Value = 5077
$mesh.Range("K12").Value = 5119

$mesh.Range("K13").Formula = "=_xlfn.VAR.S(K2:K11)"
$mesh.Range("K14").Formula = "=1-_xlfn.T.TEST(J2:J11,K2:K11,2,3)"
$mesh.Range("K15").Formula = "=J12/K12"
$mesh.Range("K16").Formula = "=B12/K12"

$mesh.Range("B15:J16").FormatConditions.Item(1).ModifyAppliesToRange($mesh.Range("B15:K16")) | Out-Null

$mesh.Range("A1").Select() | Out-Null
$mesh.Range("K4").Select() | Out-Null

$mesh.Activate() | Out-Null
